# Generate Report for Handback
# Updates the handback-status workbook: the first file's GUID changes from
# 7860255e-61b4-4612-a2bb-63b1b2ec886e to 3296156c-52af-426c-8923-1d83c7355fce,
# the second file's name changes from 81be7453-d1f2-4473-bc59-0b9cc878ac7e to
# ffffed362080-70f6-49e8-aea2-7c6384667739, and the related handback file names
# / timestamps are regenerated (and, for this run, both rows end up pointing at
# the same, newer handback xliff file + timestamp).

$wb = $excel.ActiveWorkbook

$oldGuid1 = "7860255e-61b4-4612-a2bb-63b1b2ec886e"
$newGuid1 = "3296156c-52af-426c-8923-1d83c7355fce"
$oldGuid2 = "81be7453-d1f2-4473-bc59-0b9cc878ac7e"
$newGuid2 = "ffffed362080-70f6-49e8-aea2-7c6384667739"

$newMd1 = "$newGuid1.md"
$newMd2 = "$newGuid2.md"

$newHandoffDate = "2016-08-31 03:11:51"

$newXlfZhCn = "$newGuid1.2433d031082fbb8de57c7ccd058689c8a93d32a6.zh-cn.xlf"
$newHandbackZhCnDate1 = "2016-08-31 03:11:47"
$newHandbackZhCnDate2 = "2016-08-31 03:12:18"

$newXlfDeDe = "$newGuid1.2433d031082fbb8de57c7ccd058689c8a93d32a6.de-de.xlf"
$newHandbackDeDeDate = "2016-08-31 03:12:25"

# Original (unchanged) hyperlink target addresses, keyed by worksheet + ref.
$links = @{
  "Overview" = @{
    "B2" = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5680811fc8899b4c0bd0069f4a18a689c58b2b3f/e2e/$oldGuid1.md"
    "B3" = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5680811fc8899b4c0bd0069f4a18a689c58b2b3f/e2e/$oldGuid2.md"
  }
  "zh-cn" = @{
    "A2" = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5680811fc8899b4c0bd0069f4a18a689c58b2b3f/e2e/$oldGuid1.md"
    "I2" = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c090406f33194e97219d6b2745a8e92f8ec40a9e/e2e/$oldGuid1.md"
    "A3" = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5680811fc8899b4c0bd0069f4a18a689c58b2b3f/e2e/$oldGuid2.md"
    "I3" = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c090406f33194e97219d6b2745a8e92f8ec40a9e/e2e/$oldGuid2.md"
  }
  "de-de" = @{
    "A2" = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5680811fc8899b4c0bd0069f4a18a689c58b2b3f/e2e/$oldGuid1.md"
    "I2" = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/73abe0c09d308424d35974aca99b3fb45e8955fb/e2e/$oldGuid1.md"
    "A3" = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5680811fc8899b4c0bd0069f4a18a689c58b2b3f/e2e/$oldGuid2.md"
    "I3" = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/73abe0c09d308424d35974aca99b3fb45e8955fb/e2e/$oldGuid2.md"
  }
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMd1
$wsOverview.Range("B2").Value = "e2e\$newMd1"
$wsOverview.Range("G2").Value = $newHandoffDate

$wsOverview.Range("A3").Value = $newMd2
$wsOverview.Range("B3").Value = "e2e\$newMd2"
$wsOverview.Range("G3").Value = $newHandoffDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $links["Overview"]["B2"], [type]::Missing, [type]::Missing, "e2e\$newMd1")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $links["Overview"]["B3"], [type]::Missing, [type]::Missing, "e2e\$newMd2")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newMd1
$wsZhCn.Range("G2").Value = $newXlfZhCn
$wsZhCn.Range("H2").Value = $newHandbackZhCnDate1
$wsZhCn.Range("I2").Value = $newMd1
$wsZhCn.Range("J2").Value = $newXlfZhCn
$wsZhCn.Range("K2").Value = $newHandbackZhCnDate2

$wsZhCn.Range("A3").Value = $newMd2
$wsZhCn.Range("G3").Value = $newXlfZhCn
$wsZhCn.Range("H3").Value = $newHandbackZhCnDate1
$wsZhCn.Range("I3").Value = $newMd2
$wsZhCn.Range("J3").Value = $newXlfZhCn
$wsZhCn.Range("K3").Value = $newHandbackZhCnDate2

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $links["zh-cn"]["A2"], [type]::Missing, [type]::Missing, $newMd1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $links["zh-cn"]["I2"], [type]::Missing, [type]::Missing, $newMd1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $links["zh-cn"]["A3"], [type]::Missing, [type]::Missing, $newMd2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $links["zh-cn"]["I3"], [type]::Missing, [type]::Missing, $newMd2)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newMd1
$wsDeDe.Range("G2").Value = $newXlfDeDe
$wsDeDe.Range("H2").Value = $newHandoffDate
$wsDeDe.Range("I2").Value = $newMd1
$wsDeDe.Range("J2").Value = $newXlfDeDe
$wsDeDe.Range("K2").Value = $newHandbackDeDeDate

$wsDeDe.Range("A3").Value = $newMd2
$wsDeDe.Range("G3").Value = $newXlfDeDe
$wsDeDe.Range("H3").Value = $newHandoffDate
$wsDeDe.Range("I3").Value = $newMd2
$wsDeDe.Range("J3").Value = $newXlfDeDe
$wsDeDe.Range("K3").Value = $newHandbackDeDeDate

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $links["de-de"]["A2"], [type]::Missing, [type]::Missing, $newMd1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $links["de-de"]["I2"], [type]::Missing, [type]::Missing, $newMd1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $links["de-de"]["A3"], [type]::Missing, [type]::Missing, $newMd2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $links["de-de"]["I3"], [type]::Missing, [type]::Missing, $newMd2)
